$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2..17 for columns I (I0) and J (IF)
$values = @{
    2  = @(1, 5)
    3  = @(5, 8)
    4  = @(7, 7)
    5  = @(4, 5)
    6  = @(2, 7)
    7  = @(1, 5)
    8  = @(1, 6)
    9  = @(1, 5)
    10 = @(1, 6)
    11 = @(1, 6)
    12 = @(1, 5)
    13 = @(6, 8)
    14 = @(1, 5)
    15 = @(1, 4)
    16 = @(3, 4)
    17 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
